# Daily attendance processing - 2026-01-19 14:42:33
#
# The "Recorded By" column (G) lists who/what recorded each attendance
# session. For the sessions that were (re)processed in today's run, the
# order of the two recorders is swapped from "System, <email>" to
# "<email>, System".
#
# The specific rows touched by today's run (identified by row number on
# the "Session Analysis Results" sheet):
$targetRows = @(
    2, 3, 4, 5, 6, 7,
    16, 17,
    22, 23,
    37, 38,
    43, 44,
    58, 59,
    64, 65,
    79, 80,
    85, 86, 87, 88, 89, 90,
    99, 100,
    105, 106, 107, 108, 109, 110,
    119, 120,
    125, 126, 127, 128, 129, 130,
    139, 140,
    145, 146, 147, 148, 149, 150,
    159, 160,
    165, 166, 167, 168, 169, 170,
    179, 180,
    185, 186,
    200, 201,
    206, 207,
    221, 222,
    227, 228,
    242, 243
)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
